$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The user_id_1 / 用户1 row (row 2) is being removed from the import
# template, so rows 3 and 4 (user_id_2, user_id_3) shift up to become
# rows 2 and 3.

# 1. Stash the (non-hyperlink) formatting used by the email/password
#    columns in a scratch cell so it can be restored later, since
#    re-creating hyperlinks re-applies the default hyperlink style.
$ws.Range("C2").Copy()
$ws.Range("Z1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# 2. Delete the entire row 2, shifting everything below it up by one.
$ws.Rows.Item(2).Delete()

# 3. The hyperlink objects don't move/retarget with the row shift, so
#    clear them all out and rebuild them against the new layout.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:user_id_2@1.c") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Calong@2015") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:user_id_3@1.c") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:Calong@2015") | Out-Null

# 4. Restore the original cell formatting (undo the auto hyperlink style).
$ws.Range("Z1").Copy()
$ws.Range("C2:D3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# 5. Clean up the scratch cell used for the formatting stash.
$ws.Range("Z1").Clear()

# 6. Match the saved selection from the edited workbook.
$ws.Range("F19").Select()
